# Apply the edits described by the commit:
# "introduced email check in AddAndSearchBrandNewVisitor. Also changed the tags in runner class."
#
# Observable effects on the workbook itself:
#  1. The "payment" sheet becomes the active/selected tab (was "applicant").
#  2. payment!A1 changes from a hard-coded tokenator URL to the literal label "URL".
#  3. applicant!A1 also becomes the literal label "URL", and two real tokenator
#     URLs (for chapterId=816) are now listed below it in A2 and A3.

$wb = $excel.ActiveWorkbook

$wsPayment   = $wb.Worksheets.Item("payment")
$wsApplicant = $wb.Worksheets.Item("applicant")

# payment sheet: header label only
$wsPayment.Range("A1").Value = "URL"

# applicant sheet: header label + the two new tokenator URLs
$wsApplicant.Range("A1").Value = "URL"
$wsApplicant.Range("A2").Value = "https://www.bnitest2.sasyadev.com/web/open/tokenator?concept=connect&token=f3K7zFEZcTG%2Becd7l%2BpcnVl%2FEQyClb0%2B2uuWO%2FIMB6X7lkCx8KShSd%2BGIKHzMfwPxDfikl4QLOHQC%2FHoi4NhEWGmvbOk1NlfQLPJFEPL5AEupitSDXlCZJcAlmnvOFQ4vPNW2iU7m9hGxslkmoqrNM8Q5%2FXOYnCN90YaChNd%2BwZ8dqv9YF2DuA%3D%3D&chapterId=816&step=applicant"
$wsApplicant.Range("A3").Value = "https://www.bnitest2.sasyadev.com/web/open/tokenator?concept=connect&token=ejoZQWvsSSbsNAsbhKqSEAUTa9uO%2FE4XoiGFjGYWGniD76LGo6FTYCbUFuTC64EMAQyYzCremKISOb6%2Bpn0jfsFjZ5IG2DyL4pBOv8jvblAPgZbml6Qm4AqA09V60g7s3hCRuqG1oRBcwY3dJSSPigtFSI9D9rl6n4mTmwSjoJmkDsuHZ9jChw%3D%3D&chapterId=816&step=applicant"

# Make "payment" the selected/active sheet (workbook activeTab goes from 1 -> 0)
$wsPayment.Activate()

$wb.Save()
